$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    ,@("D2", "52.332.64", $false)
    ,@("E2", "  -13.55%  ", $false)
    ,@("D3", "2.295.67", $false)
    ,@("E3", "  -20.88%  ", $false)
    ,@("E4", "  +0.06%  ", $false)
    ,@("D5", "443.30", $true)
    ,@("E5", "  -15.74%  ", $false)
    ,@("D6", "119.74", $true)
    ,@("E6", "  -15.86%  ", $false)
    ,@("D7", "0.997", $true)
    ,@("E7", "  -0.21%  ", $false)
    ,@("E8", "  -15.61%  ", $false)
    ,@("D9", "2.291.86", $false)
    ,@("E9", "  -21.18%  ", $false)
    ,@("D10", "5.26", $true)
    ,@("E10", "  -11.83%  ", $false)
    ,@("D11", "0.0870", $true)
    ,@("E11", "  -18.62%  ", $false)
    ,@("E12", "  -16.65%  ", $false)
    ,@("E13", "  -5.89%  ", $false)
    ,@("D14", "52.336.86", $false)
    ,@("E14", "  -13.54%  ", $false)
    ,@("D15", "18.69", $true)
    ,@("E15", "  -17.45%  ", $false)
    ,@("E16", "  -16.45%  ", $false)
    ,@("D17", "2.313.57", $false)
    ,@("E17", "  -20.43%  ", $false)
    ,@("D18", "3.90", $true)
    ,@("E18", "  -21.49%  ", $false)
    ,@("D19", "295.74", $true)
    ,@("E19", "  -16.03%  ", $false)
    ,@("E20", "  -23.79%  ", $false)
    ,@("D21", "0.998", $true)
    ,@("E21", "  -0.16%  ", $false)
    ,@("D22", "5.60", $true)
    ,@("E22", "  -1.76%  ", $false)
    ,@("E23", "  -22.55%  ", $false)
    ,@("E24", "  -17.51%  ", $false)
    ,@("E25", "  -20.00%  ", $false)
    ,@("E26", "  -18.03%  ", $false)
    ,@("B27", "InternetComputer(DFINITY)", $false)
    ,@("C27", "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp", $false)
    ,@("D27", "6.85", $true)
    ,@("E27", "  -12.53%  ", $false)
    ,@("B28", "USDe", $false)
    ,@("C28", "https://coinranking.com/coin/exbfr2U-0+usde-usde", $false)
    ,@("D28", "0.997", $true)
    ,@("E28", "  -0.23%  ", $false)
    ,@("D29", "0.0₃0657", $false)
    ,@("E29", "  -22.00%  ", $false)
    ,@("D30", "140.53", $true)
    ,@("E30", "  -6.51%  ", $false)
    ,@("E31", "  -14.44%  ", $false)
    ,@("E32", "  -20.52%  ", $false)
    ,@("D33", "4.70", $true)
    ,@("D34", "0.819", $true)
    ,@("E34", "  -18.19%  ", $false)
    ,@("E35", "  -21.75%  ", $false)
    ,@("D36", "0.993", $true)
    ,@("E36", "  -0.41%  ", $false)
    ,@("D37", "0.984", $true)
    ,@("E37", "  -17.77%  ", $false)
    ,@("D38", "31.67", $true)
    ,@("E38", "  -16.00%  ", $false)
    ,@("E39", "  -1.75%  ", $false)
    ,@("D40", "0.550", $true)
    ,@("E40", "  -15.05%  ", $false)
    ,@("E41", "  -13.89%  ", $false)
    ,@("E42", "  -17.05%  ", $false)
    ,@("D43", "1.904.62", $false)
    ,@("E43", "  -16.74%  ", $false)
    ,@("D44", "1.16", $true)
    ,@("E44", "  -20.84%  ", $false)
    ,@("D45", "0.0205", $true)
    ,@("E45", "  -13.63%  ", $false)
    ,@("D46", "0.0818", $true)
    ,@("E46", "  -11.19%  ", $false)
    ,@("D47", "4.16", $true)
    ,@("E47", "  -16.38%  ", $false)
    ,@("D48", "15.44", $true)
    ,@("E48", "  -24.68%  ", $false)
    ,@("E49", "  -5.06%  ", $false)
    ,@("E50", "  -13.75%  ", $false)
    ,@("D51", "14.90", $true)
    ,@("E51", "  -18.58%  ", $false)
)

foreach ($u in $updates) {
    $cellRef = $u[0]
    $newVal = $u[1]
    $forceText = $u[2]
    $rng = $ws.Range($cellRef)
    if ($forceText) {
        $rng.NumberFormat = "@"
    }
    $rng.Value = $newVal
}
